$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "table": fill in the newly-finished grid-search results for the
# "half cutoff = -1" condition (rows 10, 11, 20, 21, 22 were missing their
# Val Acc / Train Acc columns).
# ---------------------------------------------------------------------------
$wsTable = $wb.Worksheets.Item("table")

$wsTable.Range("E10").Value = 0.565
$wsTable.Range("F10").Value = 0.712

$wsTable.Range("E11").Value = 0.567
$wsTable.Range("F11").Value = 0.719

$wsTable.Range("E20").Value = 0.558
$wsTable.Range("F20").Value = 0.695

$wsTable.Range("E21").Value = 0.56
$wsTable.Range("F21").Value = 0.702

$wsTable.Range("E22").Value = 0.562
$wsTable.Range("F22").Value = 0.706

# ---------------------------------------------------------------------------
# Sheet "raw results": append the raw log lines for the 5 runs that just
# finished (same style as the rest of column A - copy formatting down).
# ---------------------------------------------------------------------------
$wsRaw = $wb.Worksheets.Item("raw results")

$separator = "=================================================================================="

$newRows = @(
    @{ Row = 69; Text = $separator },
    @{ Row = 70; Text = "pos_is_causal: False  | pos_normalize_magnitude: True  | pos_exponent: 1.5" },
    @{ Row = 71; Text = "Validation accuracy: 0.565, Training accuracy:0.712" },

    @{ Row = 73; Text = $separator },
    @{ Row = 74; Text = "pos_is_causal: False  | pos_normalize_magnitude: True  | pos_exponent: 2" },
    @{ Row = 75; Text = $separator },
    @{ Row = 76; Text = "Validation accuracy: 0.567, Training accuracy:0.719" },

    @{ Row = 78; Text = $separator },
    @{ Row = 79; Text = "pos_is_causal: True  | pos_normalize_magnitude: True  | pos_exponent: 1" },
    @{ Row = 80; Text = $separator },
    @{ Row = 81; Text = "Validation accuracy: 0.558, Training accuracy:0.695" },

    @{ Row = 83; Text = $separator },
    @{ Row = 84; Text = "pos_is_causal: True  | pos_normalize_magnitude: True  | pos_exponent: 1.5" },
    @{ Row = 85; Text = $separator },
    @{ Row = 86; Text = "Validation accuracy: 0.560, Training accuracy:0.702" },

    @{ Row = 88; Text = $separator },
    @{ Row = 89; Text = "pos_is_causal: True  | pos_normalize_magnitude: True  | pos_exponent: 2" },
    @{ Row = 90; Text = $separator },
    @{ Row = 91; Text = "Validation accuracy: 0.562, Training accuracy:0.706" }
)

# Reuse the formatting already used throughout column A (Courier New font,
# left/center aligned - style index 1) by copying it down onto the new cells.
$wsRaw.Range("A2").Copy()
foreach ($item in $newRows) {
    $cell = $wsRaw.Cells.Item($item.Row, 1)
    $cell.PasteSpecial(-4122)
    if ($item.Text -eq $separator) {
        # The separator line reuses the very first shared-string entry
        # (index 0); assigning it through .Value on several cells confuses
        # the engine's dirty-tracking, so paste the value across instead.
        $cell.PasteSpecial(-4163)
    } else {
        $cell.Value = $item.Text
    }
}
$excel.CutCopyMode = $false

# Update selections/scroll position to mirror where the author ended up
# after finishing the edits. "table" must remain the active/selected tab,
# so touch "raw results" first and re-activate "table" last.
$wsRaw.Activate()
$wsRaw.Range("A93").Select()
$excel.ActiveWindow.ScrollRow = 65
$excel.ActiveWindow.ScrollColumn = 1

$wsTable.Activate()
$wsTable.Range("F23").Select()
